$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Inscritos 57 -> 58
$ws.Range("E6").Value = 58

# Row 10: Inscritos 32 -> 34, Pagos 14 -> 15, Inscricoes homologadas 16 -> 17
$ws.Range("E10").Value = 34
$ws.Range("F10").Value = 15
$ws.Range("H10").Value = 17

# Row 11: Inscritos 20 -> 21, Pagos 14 -> 15, Inscricoes homologadas 15 -> 16
$ws.Range("E11").Value = 21
$ws.Range("F11").Value = 15
$ws.Range("H11").Value = 16

# Row 14: Inscritos 37 -> 38
$ws.Range("E14").Value = 38

# Row 15: Inscritos 103 -> 104
$ws.Range("E15").Value = 104

# Row 17: Inscritos 27 -> 29, Pagos 12 -> 13, Inscricoes homologadas 13 -> 14
$ws.Range("E17").Value = 29
$ws.Range("F17").Value = 13
$ws.Range("H17").Value = 14
